# Update the cryptocurrency price/volume table with freshly scraped values.
# Columns D (Price) and E (Volume(1h)) hold number-like text (e.g. "0.06630",
# "1.826.24") that Excel would otherwise auto-convert to a Number (losing
# formatting like trailing zeros or multiple "." thousands separators), so
# those values are entered with a leading apostrophe to force literal text,
# then the style is reset to Normal so no stray quote-prefix indicator/format
# is left behind on the cell (matches the original plain inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '27.466.27' }
    @{ Cell = "E2"; Value = '  -0.40%  ' }
    @{ Cell = "D3"; Value = '1.825.96' }
    @{ Cell = "E3"; Value = '  -1.76%  ' }
    @{ Cell = "D4"; Value = '1.006' }
    @{ Cell = "E4"; Value = '  -0.72%  ' }
    @{ Cell = "D5"; Value = '331.35' }
    @{ Cell = "E5"; Value = '  -0.54%  ' }
    @{ Cell = "E6"; Value = '  -0.70%  ' }
    @{ Cell = "D7"; Value = '0.4557' }
    @{ Cell = "E7"; Value = '  -1.97%  ' }
    @{ Cell = "D8"; Value = '0.3805' }
    @{ Cell = "E8"; Value = '  -1.96%  ' }
    @{ Cell = "D9"; Value = '46.39' }
    @{ Cell = "E9"; Value = '  +1.32%  ' }
    @{ Cell = "D10"; Value = '0.07878' }
    @{ Cell = "E10"; Value = '  -0.87%  ' }
    @{ Cell = "D11"; Value = '0.9714' }
    @{ Cell = "E11"; Value = '  -2.20%  ' }
    @{ Cell = "D12"; Value = '20.99' }
    @{ Cell = "E12"; Value = '  -2.16%  ' }
    @{ Cell = "D13"; Value = '1.826.24' }
    @{ Cell = "E13"; Value = '  -1.72%  ' }
    @{ Cell = "D14"; Value = '5.865' }
    @{ Cell = "E14"; Value = '  -1.49%  ' }
    @{ Cell = "D15"; Value = '7.029' }
    @{ Cell = "E15"; Value = '  -1.89%  ' }
    @{ Cell = "D16"; Value = '1.007' }
    @{ Cell = "E16"; Value = '  -0.71%  ' }
    @{ Cell = "D17"; Value = '89.15' }
    @{ Cell = "E17"; Value = '  +1.85%  ' }
    @{ Cell = "D18"; Value = '0.06630' }
    @{ Cell = "E18"; Value = '  -1.09%  ' }
    @{ Cell = "E19"; Value = '  -1.20%  ' }
    @{ Cell = "D20"; Value = '17.17' }
    @{ Cell = "E20"; Value = '  +2.04%  ' }
    @{ Cell = "E21"; Value = '  -0.66%  ' }
    @{ Cell = "D22"; Value = '27.449.10' }
    @{ Cell = "E22"; Value = '  -0.51%  ' }
    @{ Cell = "D23"; Value = '5.326' }
    @{ Cell = "E23"; Value = '  -1.94%  ' }
    @{ Cell = "D24"; Value = '10.79' }
    @{ Cell = "E24"; Value = '  -0.15%  ' }
    @{ Cell = "E25"; Value = '  -0.49%  ' }
    @{ Cell = "D26"; Value = '2.042.67' }
    @{ Cell = "E26"; Value = '  -1.88%  ' }
    @{ Cell = "D27"; Value = '156.78' }
    @{ Cell = "E27"; Value = '  -1.07%  ' }
    @{ Cell = "D28"; Value = '19.40' }
    @{ Cell = "E28"; Value = '  -1.19%  ' }
    @{ Cell = "D29"; Value = '2.060' }
    @{ Cell = "E29"; Value = '  -1.32%  ' }
    @{ Cell = "D30"; Value = '5.247' }
    @{ Cell = "E30"; Value = '  -1.58%  ' }
    @{ Cell = "D31"; Value = '118.13' }
    @{ Cell = "E31"; Value = '  -2.45%  ' }
    @{ Cell = "D32"; Value = '0.9467' }
    @{ Cell = "E32"; Value = '  -1.85%  ' }
    @{ Cell = "D33"; Value = '0.09308' }
    @{ Cell = "E33"; Value = '  -1.30%  ' }
    @{ Cell = "D34"; Value = '3.580' }
    @{ Cell = "E34"; Value = '  -1.73%  ' }
    @{ Cell = "D35"; Value = '5.235' }
    @{ Cell = "E35"; Value = '  -0.57%  ' }
    @{ Cell = "D36"; Value = '1.321' }
    @{ Cell = "E36"; Value = '  +0.77%  ' }
    @{ Cell = "D37"; Value = '0.05922' }
    @{ Cell = "E37"; Value = '  -1.42%  ' }
    @{ Cell = "D38"; Value = '0.02179' }
    @{ Cell = "E38"; Value = '  -1.33%  ' }
    @{ Cell = "D39"; Value = '1.161' }
    @{ Cell = "E39"; Value = '  -2.89%  ' }
    @{ Cell = "D40"; Value = '8.021' }
    @{ Cell = "E40"; Value = '  -0.95%  ' }
    @{ Cell = "D41"; Value = '0.5742' }
    @{ Cell = "E41"; Value = '  -2.22%  ' }
    @{ Cell = "D42"; Value = '0.1829' }
    @{ Cell = "E42"; Value = '  -2.19%  ' }
    @{ Cell = "D43"; Value = '10.01' }
    @{ Cell = "E43"; Value = '  -1.33%  ' }
    @{ Cell = "D44"; Value = '1.268' }
    @{ Cell = "E44"; Value = '  +1.13%  ' }
    @{ Cell = "B45"; Value = 'Decentraland' }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = "D45"; Value = '0.5444' }
    @{ Cell = "E45"; Value = '  -2.58%  ' }
    @{ Cell = "B46"; Value = 'EnergySwap' }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = "D46"; Value = '11.94' }
    @{ Cell = "E46"; Value = '  -0.48%  ' }
    @{ Cell = "D47"; Value = '1.863' }
    @{ Cell = "E47"; Value = '  -2.03%  ' }
    @{ Cell = "D48"; Value = '0.06611' }
    @{ Cell = "E48"; Value = '  -2.04%  ' }
    @{ Cell = "D49"; Value = '110.31' }
    @{ Cell = "E49"; Value = '  -1.49%  ' }
    @{ Cell = "D50"; Value = '1.039' }
    @{ Cell = "E50"; Value = '  -1.55%  ' }
    @{ Cell = "E51"; Value = '  -0.75%  ' }
)

foreach ($u in $updates) {
    $col = $u.Cell.Substring(0, 1)
    $value = $u.Value
    $range = $ws.Range($u.Cell)
    if ($col -eq 'D' -or $col -eq 'E') {
        # Force text (prevents Excel turning e.g. '0.06630' or '1.006' into a Number)
        $range.Value = "'" + $value
        # Drop the quote-prefix style Excel applies for apostrophe-led entries so
        # the cell keeps its original (default) style, same as the source file.
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

